# Week 4 Lecture 2
#
# Slide 35 ("Lecture Recap" / for-loop recap slide) is turned into a big,
# centered "PRACTICE!" title slide: the Content Placeholder and the small
# "Practice!" textbox are removed, and the Title placeholder is repositioned,
# resized, given a fixed (no-autofit) body and a big centered "PRACTICE!"
# run (with the "!" in a light purple).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(35)

# Drop the body content placeholder and the small "Practice!" textbox
# entirely (Cut fully removes the shape from the slide, unlike Delete on a
# placeholder which just resets its contents).
$s.Shapes.Item(2).Cut()
$s.Shapes.Item(2).Cut()

# Reposition / resize the remaining Title placeholder and turn autofit off.
$title = $s.Shapes.Item(1)
$title.Left = 163.76925659842522
$title.Top = 203.62893673779527
$title.Width = 632.4616089031496
$title.Height = 132.742164634252
$title.TextFrame.AutoSize = 0

# Replace the title text with a big, bold, centered "PRACTICE!" where the
# exclamation point is highlighted in a light purple.
$tr = $title.TextFrame.TextRange
$tr.Text = "PRACTICE!"
$tr.Font.Size = 120
$tr.Font.Bold = $true
$tr.ParagraphFormat.Alignment = 2

$bang = $tr.Characters(9, 1)
$bang.Font.Color.RGB = 16751052
